# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    # Force text storage so numeric-looking strings (e.g. "515.85",
    # "0.104") are not silently reinterpreted as numbers, while
    # keeping the cell on the default "Normal" style (no custom s=).
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "55.902.06"
$ws.Range("E2").Value = "  +3.77%  "
Set-TextCell $ws "D3" "2.306.06"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextCell $ws "D5" "515.85"
$ws.Range("E5").Value = "  +4.63%  "
Set-TextCell $ws "D6" "132.38"
$ws.Range("E6").Value = "  +3.49%  "
Set-TextCell $ws "D7" "0.993"
$ws.Range("E7").Value = "  -0.68%  "
Set-TextCell $ws "D8" "0.534"
$ws.Range("E8").Value = "  +2.15%  "
Set-TextCell $ws "D9" "2.329.13"
$ws.Range("E9").Value = "  +3.28%  "
Set-TextCell $ws "D10" "0.104"
$ws.Range("E10").Value = "  +10.25%  "
$ws.Range("E11").Value = "  +1.30%  "
Set-TextCell $ws "D12" "5.11"
$ws.Range("E12").Value = "  +8.36%  "
Set-TextCell $ws "D13" "0.344"
$ws.Range("E13").Value = "  +3.03%  "
Set-TextCell $ws "D14" "24.18"
$ws.Range("E14").Value = "  +7.10%  "
Set-TextCell $ws "D15" "2.717.48"
$ws.Range("E15").Value = "  +2.44%  "
Set-TextCell $ws "D16" "56.008.82"
$ws.Range("E16").Value = "  +4.00%  "
Set-TextCell $ws "D17" "0.0000136"
$ws.Range("E17").Value = "  +5.50%  "
Set-TextCell $ws "D18" "2.345.84"
$ws.Range("E18").Value = "  +3.33%  "
Set-TextCell $ws "D19" "10.57"
$ws.Range("E19").Value = "  +3.75%  "
Set-TextCell $ws "D20" "4.25"
$ws.Range("E20").Value = "  +3.21%  "
Set-TextCell $ws "D21" "320.67"
$ws.Range("E21").Value = "  +6.65%  "
Set-TextCell $ws "D22" "6.65"
$ws.Range("E22").Value = "  +5.85%  "
Set-TextCell $ws "D23" "0.996"
$ws.Range("E23").Value = "  -0.27%  "
Set-TextCell $ws "D24" "60.49"
$ws.Range("E24").Value = "  -0.17%  "
Set-TextCell $ws "D25" "0.991"
$ws.Range("E25").Value = "  -0.73%  "
Set-TextCell $ws "D26" "0.158"
$ws.Range("E26").Value = "  +7.04%  "
Set-TextCell $ws "D27" "7.67"
$ws.Range("E27").Value = "  +5.73%  "
Set-TextCell $ws "D28" "171.93"
$ws.Range("E28").Value = "  +1.09%  "
Set-TextCell $ws "D29" "1.19"
$ws.Range("E29").Value = "  +10.36%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws "D30" "0.0₃0726"
$ws.Range("E30").Value = "  +6.38%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D31" "6.25"
$ws.Range("E31").Value = "  +5.84%  "
Set-TextCell $ws "D32" "1.68"
$ws.Range("E32").Value = "  +5.36%  "
Set-TextCell $ws "D33" "18.26"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("E34").Value = "  -0.04%  "
Set-TextCell $ws "D35" "0.990"
$ws.Range("E35").Value = "  -0.89%  "
Set-TextCell $ws "D36" "1.26"
$ws.Range("E36").Value = "  +6.56%  "
Set-TextCell $ws "D37" "0.924"
$ws.Range("E37").Value = "  -1.25%  "
Set-TextCell $ws "D38" "3.99"
$ws.Range("E38").Value = "  +8.26%  "
Set-TextCell $ws "D39" "1.52"
$ws.Range("E39").Value = "  +9.84%  "
Set-TextCell $ws "D40" "37.29"
$ws.Range("E40").Value = "  +4.13%  "
Set-TextCell $ws "D41" "0.384"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D42" "3.60"
$ws.Range("E42").Value = "  +7.92%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws "D43" "138.68"
$ws.Range("E43").Value = "  +10.58%  "
Set-TextCell $ws "D44" "5.14"
$ws.Range("E44").Value = "  +8.16%  "
Set-TextCell $ws "D45" "266.53"
$ws.Range("E45").Value = "  +12.00%  "
Set-TextCell $ws "D46" "0.0511"
$ws.Range("E46").Value = "  +4.69%  "
Set-TextCell $ws "D47" "0.0927"
$ws.Range("E47").Value = "  +4.48%  "
Set-TextCell $ws "D48" "0.556"
$ws.Range("E48").Value = "  +3.00%  "
Set-TextCell $ws "D49" "0.384"
$ws.Range("E49").Value = "  +4.19%  "
Set-TextCell $ws "D50" "0.0216"
$ws.Range("E50").Value = "  +6.49%  "
Set-TextCell $ws "D51" "16.89"
$ws.Range("E51").Value = "  +5.52%  "
